# Integrated report generation and appending to report in main()
#
# This script appends a freshly generated analysis report (timestamp
# 14:32:42 21/01/2026) to the SOLUSDT_analysis workbook:
#   - refreshes the "Dashboard" summary cells with the latest values
#   - appends a new data row to "Data" (row 22) and to each of the
#     per-topic sheets (row 24): "Price Analysis", "Technical Analysis",
#     "Fundamental Analysis", "Sentiment Analysis", "Predictions"
#   - extends every chart series range so the new row is plotted

$wb = $excel.ActiveWorkbook

# Helper: write a value that must stay a literal text cell even when it
# looks numeric/currency/percentage (e.g. "$127.21", "-1.30%"), avoiding
# Excel's usual smart "looks like a number" auto-conversion on .Value
# assignment, and without leaving a residual non-default cell style behind.
function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# Helper: write the shared "datetime" column (numeric date serial with the
# yyyy-mm-dd h:mm:ss display format used throughout the workbook).
function Set-DateCell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
    $ws.Range($addr).NumberFormat = "yyyy-mm-dd h:mm:ss"
}

# ---------------------------------------------------------------------
# 1. Dashboard - "Latest Analysis" + "Summary Statistics" refresh
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")
Set-TextValue $dash "B4" "14:32:42 21/01/2026"
Set-TextValue $dash "B5" "$127.21"
Set-TextValue $dash "B6" "-1.30%"
Set-TextValue $dash "B7" "$122.57"
Set-TextValue $dash "B9" "0.013"
Set-TextValue $dash "B10" "0.57"
$dash.Range("B15").Value = 21
Set-TextValue $dash "B16" "$127.50"
Set-TextValue $dash "B18" "12.19%"
Set-TextValue $dash "B19" "32.40"
Set-TextValue $dash "B20" "0.027"

# ---------------------------------------------------------------------
# 2. Data sheet - append the new raw report row (row 22)
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("Data")
Set-TextValue $data "A22" "14:32:42 21/01/2026"
Set-DateCell  $data "B22" 46043.60604952029
$data.Range("C22").Value = 127.21
$data.Range("D22").Value = 129.77
$data.Range("E22").Value = 124.68
$data.Range("F22").Value = -1.67
$data.Range("G22").Value = -1.296
$data.Range("H22").Value = 2697859.957
$data.Range("I22").Value = 344386807.26149
$data.Range("J22").Value = 1356132
$data.Range("K22").Value = 4.5
$data.Range("L22").Value = 5.5
$data.Range("M22").Value = 4.5
$data.Range("N22").Value = 7
$data.Range("O22").Value = 2
Set-TextValue $data "P22" "bearish"
$data.Range("Q22").Value = 1
$data.Range("R22").Value = -0.01
Set-TextValue $data "S22" "low"
$data.Range("T22").Value = 49.73
Set-TextValue $data "U22" "neutral"
$data.Range("V22").Value = -0.8659
Set-TextValue $data "W22" "bullish"
$data.Range("X22").Value = 0.013
$data.Range("Y22").Value = 0.022
$data.Range("Z22").Value = 0
$data.Range("AA22").Value = 58.2
Set-TextValue $data "AB22" "Neutral"
$data.Range("AC22").Value = 0.451
$data.Range("AD22").Value = 0.549
$data.Range("AE22").Value = 9.699999999999999
$data.Range("AF22").Value = 122.57
$data.Range("AG22").Value = -3.65
Set-TextValue $data "AH22" "downward"
$data.Range("AI22").Value = 126.4
$data.Range("AJ22").Value = 127.79
$data.Range("AK22").Value = 126.4
$data.Range("AL22").Value = 0.17
Set-TextValue $data "AM22" "Hold or wait for clearer signals"
$data.Range("AN22").Value = 0.57

# ---------------------------------------------------------------------
# 3. Price Analysis - append row 24
# ---------------------------------------------------------------------
$price = $wb.Worksheets.Item("Price Analysis")
Set-TextValue $price "A24" "14:32:42 21/01/2026"
Set-DateCell  $price "B24" 46043.60604952029
$price.Range("C24").Value = 127.21
$price.Range("D24").Value = 129.77
$price.Range("E24").Value = 124.68
$price.Range("F24").Value = -1.67
$price.Range("G24").Value = -1.296
$price.Range("H24").Value = 2697859.957

# ---------------------------------------------------------------------
# 4. Technical Analysis - append row 24
# ---------------------------------------------------------------------
$tech = $wb.Worksheets.Item("Technical Analysis")
Set-TextValue $tech "A24" "14:32:42 21/01/2026"
Set-DateCell  $tech "B24" 46043.60604952029
$tech.Range("C24").Value = 1
$tech.Range("D24").Value = -0.01
$tech.Range("E24").Value = 49.73
$tech.Range("F24").Value = -0.8659

# ---------------------------------------------------------------------
# 5. Fundamental Analysis - append row 24
# ---------------------------------------------------------------------
$fund = $wb.Worksheets.Item("Fundamental Analysis")
Set-TextValue $fund "A24" "14:32:42 21/01/2026"
Set-DateCell  $fund "B24" 46043.60604952029
$fund.Range("C24").Value = 4.5
$fund.Range("D24").Value = 5.5
$fund.Range("E24").Value = 4.5
$fund.Range("F24").Value = 7
$fund.Range("G24").Value = 2

# ---------------------------------------------------------------------
# 6. Sentiment Analysis - append row 24
# ---------------------------------------------------------------------
$sent = $wb.Worksheets.Item("Sentiment Analysis")
Set-TextValue $sent "A24" "14:32:42 21/01/2026"
Set-DateCell  $sent "B24" 46043.60604952029
$sent.Range("C24").Value = 0.013
$sent.Range("D24").Value = 0.022
$sent.Range("E24").Value = 0
$sent.Range("F24").Value = 58.2

# ---------------------------------------------------------------------
# 7. Predictions - append row 24
# ---------------------------------------------------------------------
$pred = $wb.Worksheets.Item("Predictions")
Set-TextValue $pred "A24" "14:32:42 21/01/2026"
Set-DateCell  $pred "B24" 46043.60604952029
$pred.Range("C24").Value = 127.21
$pred.Range("D24").Value = 122.57
$pred.Range("E24").Value = -3.65
$pred.Range("F24").Value = 0.451
$pred.Range("G24").Value = 0.549
$pred.Range("H24").Value = 126.4
$pred.Range("I24").Value = 127.79

# ---------------------------------------------------------------------
# 8. Extend every chart series so the new row is plotted ($4:$23 -> $4:$24)
# ---------------------------------------------------------------------

# Price Analysis charts (drawing1): chart1 (3 series: C,D,E), chart2 (1 series: H)
$co = $price.ChartObjects().Item(1).Chart
$co.SeriesCollection(1).Values = "='Price Analysis'!`$C`$4:`$C`$24"
$co.SeriesCollection(1).XValues = "='Price Analysis'!`$B`$4:`$B`$24"
$co.SeriesCollection(2).Values = "='Price Analysis'!`$D`$4:`$D`$24"
$co.SeriesCollection(2).XValues = "='Price Analysis'!`$B`$4:`$B`$24"
$co.SeriesCollection(3).Values = "='Price Analysis'!`$E`$4:`$E`$24"
$co.SeriesCollection(3).XValues = "='Price Analysis'!`$B`$4:`$B`$24"

$co = $price.ChartObjects().Item(2).Chart
$co.SeriesCollection(1).Values = "='Price Analysis'!`$H`$4:`$H`$24"
$co.SeriesCollection(1).XValues = "='Price Analysis'!`$B`$4:`$B`$24"

# Technical Analysis charts (drawing2): chart3 (1 series: E), chart4 (2 series: C,D)
$co = $tech.ChartObjects().Item(1).Chart
$co.SeriesCollection(1).Values = "='Technical Analysis'!`$E`$4:`$E`$24"
$co.SeriesCollection(1).XValues = "='Technical Analysis'!`$B`$4:`$B`$24"

$co = $tech.ChartObjects().Item(2).Chart
$co.SeriesCollection(1).Values = "='Technical Analysis'!`$C`$4:`$C`$24"
$co.SeriesCollection(1).XValues = "='Technical Analysis'!`$B`$4:`$B`$24"
$co.SeriesCollection(2).Values = "='Technical Analysis'!`$D`$4:`$D`$24"
$co.SeriesCollection(2).XValues = "='Technical Analysis'!`$B`$4:`$B`$24"

# Fundamental Analysis charts (drawing3): chart5 (5 series: C,D,E,F,G)
$co = $fund.ChartObjects().Item(1).Chart
$co.SeriesCollection(1).Values = "='Fundamental Analysis'!`$C`$4:`$C`$24"
$co.SeriesCollection(1).XValues = "='Fundamental Analysis'!`$B`$4:`$B`$24"
$co.SeriesCollection(2).Values = "='Fundamental Analysis'!`$D`$4:`$D`$24"
$co.SeriesCollection(2).XValues = "='Fundamental Analysis'!`$B`$4:`$B`$24"
$co.SeriesCollection(3).Values = "='Fundamental Analysis'!`$E`$4:`$E`$24"
$co.SeriesCollection(3).XValues = "='Fundamental Analysis'!`$B`$4:`$B`$24"
$co.SeriesCollection(4).Values = "='Fundamental Analysis'!`$F`$4:`$F`$24"
$co.SeriesCollection(4).XValues = "='Fundamental Analysis'!`$B`$4:`$B`$24"
$co.SeriesCollection(5).Values = "='Fundamental Analysis'!`$G`$4:`$G`$24"
$co.SeriesCollection(5).XValues = "='Fundamental Analysis'!`$B`$4:`$B`$24"

# Sentiment Analysis charts (drawing4): chart6 (3 series: C,D,E), chart7 (1 series: F)
$co = $sent.ChartObjects().Item(1).Chart
$co.SeriesCollection(1).Values = "='Sentiment Analysis'!`$C`$4:`$C`$24"
$co.SeriesCollection(1).XValues = "='Sentiment Analysis'!`$B`$4:`$B`$24"
$co.SeriesCollection(2).Values = "='Sentiment Analysis'!`$D`$4:`$D`$24"
$co.SeriesCollection(2).XValues = "='Sentiment Analysis'!`$B`$4:`$B`$24"
$co.SeriesCollection(3).Values = "='Sentiment Analysis'!`$E`$4:`$E`$24"
$co.SeriesCollection(3).XValues = "='Sentiment Analysis'!`$B`$4:`$B`$24"

$co = $sent.ChartObjects().Item(2).Chart
$co.SeriesCollection(1).Values = "='Sentiment Analysis'!`$F`$4:`$F`$24"
$co.SeriesCollection(1).XValues = "='Sentiment Analysis'!`$B`$4:`$B`$24"

# Predictions charts (drawing5): chart8 (2 series: C,D), chart9 (3 series: C,H,I)
$co = $pred.ChartObjects().Item(1).Chart
$co.SeriesCollection(1).Values = "=Predictions!`$C`$4:`$C`$24"
$co.SeriesCollection(1).XValues = "=Predictions!`$B`$4:`$B`$24"
$co.SeriesCollection(2).Values = "=Predictions!`$D`$4:`$D`$24"
$co.SeriesCollection(2).XValues = "=Predictions!`$B`$4:`$B`$24"

$co = $pred.ChartObjects().Item(2).Chart
$co.SeriesCollection(1).Values = "=Predictions!`$C`$4:`$C`$24"
$co.SeriesCollection(1).XValues = "=Predictions!`$B`$4:`$B`$24"
$co.SeriesCollection(2).Values = "=Predictions!`$H`$4:`$H`$24"
$co.SeriesCollection(2).XValues = "=Predictions!`$B`$4:`$B`$24"
$co.SeriesCollection(3).Values = "=Predictions!`$I`$4:`$I`$24"
$co.SeriesCollection(3).XValues = "=Predictions!`$B`$4:`$B`$24"
